$d = $word.ActiveDocument

function Get-ParagraphXmlHeader {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
}

function Get-ParagraphXmlFooter {
    return '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- 1) First Reviewer-1 paragraph that was ONLY the "GitHub issue #136" hyperlink. ---
# It becomes a plain bold run (hyperlink removed) describing the new summary figure.
$targetText1 = "GitHub issue #136"
$paraIndex1 = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $targetText1) {
        $paraIndex1 = $i
        break
    }
}
if ($paraIndex1 -eq -1) {
    throw "Could not find the standalone 'GitHub issue #136' paragraph"
}

$p1 = $d.Paragraphs.Item($paraIndex1)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$xml1 = (Get-ParagraphXmlHeader) + '<w:p><w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">We have added a summary figure (new Fig. 3) that we believe provides a very helpful summary. Thanks for the suggestion!</w:t></w:r></w:p>' + (Get-ParagraphXmlFooter)
$r1.InsertXML($xml1)

# --- 2) The "We agree that forests are complex..." paragraph that embedded the same ---
#        hyperlink mid-sentence; merge into a single bold run & fix "anaalysis" typo.
$needle2 = "We agree that forests are complex"
$paraIndex2 = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith($needle2)) {
        $paraIndex2 = $i
        break
    }
}
if ($paraIndex2 -eq -1) {
    throw "Could not find the 'We agree that forests are complex' paragraph"
}

$p2 = $d.Paragraphs.Item($paraIndex2)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$newText2 = "We agree that forests are complex, and this method both reveals and allows us to start to distill some of that complexity. We have added a summary figure (new Fig. 3) to help to highlight some of the general patterns. There is obviously a lot more work to be done on this theme that is beyond the scope of the current analysis."
$xml2 = (Get-ParagraphXmlHeader) + '<w:p><w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">' + $newText2 + '</w:t></w:r></w:p>' + (Get-ParagraphXmlFooter)
$r2.InsertXML($xml2)

Write-Output "done"
